# Quantities.xlsx daily roll-forward update:
#   - Each existing data row (2..29) takes on the values that used to
#     belong to the next row down (the whole table shifts "up" by one day),
#     except the last existing row (29) which keeps its own values.
#   - A brand-new row (30) is appended, duplicating the former last row's
#     values, and labelled with the next calendar date.
#   - Column A holds the date label for each row; after the shift every
#     row's date moves forward by one day, and the new row gets the next
#     date after the previous last one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 29
$newRow = $lastDataRow + 1
$lastCol = 10   # column J

# 1) Snapshot the current B:J values for rows (firstDataRow+1) .. lastDataRow
#    before anything is overwritten.
$snapshot = @{}
for ($r = $firstDataRow + 1; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    for ($c = 2; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}
# also remember the (unchanged) last row, to duplicate into the new row
$lastRowVals = @{}
for ($c = 2; $c -le $lastCol; $c++) {
    $lastRowVals[$c] = $ws.Cells.Item($lastDataRow, $c).Value2
}

# 2) Shift B:J up by one row: row r gets what row r+1 used to hold.
for ($r = $firstDataRow; $r -le $lastDataRow - 1; $r++) {
    $src = $snapshot[$r + 1]
    for ($c = 2; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $src[$c]
    }
}
# row $lastDataRow (29) keeps its own values - nothing to do there.

# 3) Append the new row, duplicating the former last row's B:J values.
for ($c = 2; $c -le $lastCol; $c++) {
    $ws.Cells.Item($newRow, $c).Value = $lastRowVals[$c]
}

# 4) Roll every date label in column A forward by one day.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $oldDate = [DateTime]::ParseExact($ws.Cells.Item($r, 1).Text, "yyyy-MM-dd", $null)
    $newDate = $oldDate.AddDays(1)
    $ws.Cells.Item($r, 1).Value = "'" + $newDate.ToString("yyyy-MM-dd")
}

# 5) New row's date label: one day after the (already-shifted) previous row.
$ws.Range("A29").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$prevDate = [DateTime]::ParseExact($ws.Cells.Item($lastDataRow, 1).Text, "yyyy-MM-dd", $null)
$nextDate = $prevDate.AddDays(1)
$ws.Cells.Item($newRow, 1).Value = "'" + $nextDate.ToString("yyyy-MM-dd")

$ws.Range("A1:J" + $newRow).Select()
